# Updated fitting parameters (r_s_star, h_p_star) on the Parameters sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Parameters")
$ws.Range("J2").Value = 0.03447
$ws.Range("K2").Value = 0.0175
